# Update: UAE Premier League base update (30-03-2024)
# Swaps the row-pairs that were reordered in the source refresh, and removes
# the three not-yet-played fixtures that dropped out of the feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row pairs whose full records (everything except the running id-counter in
# column A and the Date in column E) traded places.
$pairs = @(
    @(2, 3),
    @(16, 17),
    @(18, 19),
    @(25, 26),
    @(39, 40),
    @(70, 71),
    @(98, 99),
    @(107, 108)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $bd1 = $ws.Range("B$r1`:D$r1").Value2
    $bd2 = $ws.Range("B$r2`:D$r2").Value2
    $fac1 = $ws.Range("F$r1`:AC$r1").Value2
    $fac2 = $ws.Range("F$r2`:AC$r2").Value2

    $ws.Range("B$r1`:D$r1").Value2 = $bd2
    $ws.Range("B$r2`:D$r2").Value2 = $bd1
    $ws.Range("F$r1`:AC$r1").Value2 = $fac2
    $ws.Range("F$r2`:AC$r2").Value2 = $fac1
}

# Remove the three unplayed fixtures (rows 112-114) that were dropped.
$ws.Range("A112:A114").EntireRow.Delete()
